$wb = $excel.ActiveWorkbook

# --- Update metadata sheet: Date value ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2023-05-30T18:28:01+00:00"

# --- Update Concepts sheet ---
$ws = $wb.Worksheets.Item("Concepts")

# Fix display text for TRATU row (row 12): Transcriptome -> Tumoral Transcriptome
$ws.Range("C12").Value = "Tumoral Transcriptome"

# Add a new concept row (row 15): EXTUM / Tumoral Exome
# Copy formatting (and the "Level" = 1 text value) down from the last data row (14)
$ws.Range("A14:D14").Copy()
$ws.Range("A15:D15").PasteSpecial(-4122)

# Preserve "1" as a text value (matching the rest of column A) by copying value+format
$ws.Range("A14").Copy()
$ws.Range("A15").PasteSpecial(-4163)

$ws.Range("B15").Value = "EXTUM"
$ws.Range("C15").Value = "Tumoral Exome"

$ws.Application.CutCopyMode = $false
